# Generate Report for Handoff
# Updates the "b.md" row to reflect that file "b" is now ready for handoff,
# with a freshly generated handoff package (xlf file) and timestamp, across
# the Overview sheet and each locale sheet (zh-cn, de-de).

$wb = $excel.ActiveWorkbook

$status = "Ready for handoff"

# ---------------------------------------------------------------------------
# Overview sheet: row 3 corresponds to "b.md"
# ---------------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = $status
$overview.Range("C3").Value = $status
$overview.Range("D3").Value = "2016-03-22 20:36:22"

# ---------------------------------------------------------------------------
# zh-cn sheet: row 3 corresponds to "b.md"
# ---------------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcnFile = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"

$zhcn.Range("C3").Value = $status
$zhcn.Range("D3").Value = $zhcnFile
$zhcn.Range("E3").Value = "2016-03-22 20:36:16"

foreach ($hl in $zhcn.Hyperlinks) {
    if ($hl.Range.Address() -eq '$D$3') {
        $hl.TextToDisplay = $zhcnFile
    }
}

# ---------------------------------------------------------------------------
# de-de sheet: row 3 corresponds to "b.md"
# ---------------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dedeFile = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"

$dede.Range("C3").Value = $status
$dede.Range("D3").Value = $dedeFile
$dede.Range("E3").Value = "2016-03-22 20:36:22"

foreach ($hl in $dede.Hyperlinks) {
    if ($hl.Range.Address() -eq '$D$3') {
        $hl.TextToDisplay = $dedeFile
    }
}
